# Apply the edit described by the diff:
#  - Clear A6:D6 (was OutilsFR / Liens : / Easy-It / https://easy-it.kiabi.fr/)
#  - Fill B7 with "Liens :" (keeps A7/C7/D7 = OutilsFR / Watchdoc WFR616 / url)
#  - Clear A8:D9 (was Ressources HumainesFR block)
#  - Clear A42:D42 (was DecisionnelFR / Liens : / Kiperf / https://kiperf.kiabi.pro/)
#  - Fill B43 with "Liens :" (keeps A43/C43/D43/E43 = DecisionnelFR / Zone de lancement BI / url / ignore)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pick up the correctly vertically-centered format from a neighbouring cell
# on the same row (C7 / C43) before writing the new "Liens :" labels, so the
# resulting cell style matches the other "Liens :" cells in the sheet.
$ws.Range("C7").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C43").Copy()
$ws.Range("B43").PasteSpecial(-4122)

# Now clear the rows/cells that become blank.
$ws.Range("A6:D6").ClearContents()
$ws.Range("A8:D9").ClearContents()
$ws.Range("A42:D42").ClearContents()

# Fill in the two labels that remain after the clears.
$ws.Range("B7").Value = "Liens :"
$ws.Range("B43").Value = "Liens :"

# Restore the sheet view state (top-left cell + selection) to match the
# target workbook.
$ws.Range("D42").Select()
$excel.ActiveWindow.ScrollRow = 10
